$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 114 (shift existing rows 114-173 down to 117-176)
$ws.Rows.Item(114).Insert()
$ws.Rows.Item(114).Insert()
$ws.Rows.Item(114).Insert()

# Populate the three newly inserted rows with the new "Palta" price records
$ws.Range("A114").Value = 1
$ws.Range("B114").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C114").Value = "Arica y Parinacota"
$ws.Range("D114").Value = 45016
$ws.Range("E114").Value = 15
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100106
$ws.Range("H114").Value = "Oleaginosos"
$ws.Range("I114").Value = 100106002
$ws.Range("J114").Value = "Palta"
$ws.Range("K114").Value = "Hass"
$ws.Range("L114").Value = "Cuarta"
$ws.Range("M114").Value = 650
$ws.Range("N114").Value = 20000
$ws.Range("O114").Value = 21000
$ws.Range("P114").Value = 20692
$ws.Range("Q114").Value = "`$/bandeja 10 kilos"
$ws.Range("R114").Value = "Perú"
$ws.Range("S114").Value = 2069
$ws.Range("T114").Value = 10
$ws.Range("A115").Value = 1
$ws.Range("B115").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C115").Value = "Arica y Parinacota"
$ws.Range("D115").Value = 45016
$ws.Range("E115").Value = 15
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100106
$ws.Range("H115").Value = "Oleaginosos"
$ws.Range("I115").Value = 100106002
$ws.Range("J115").Value = "Palta"
$ws.Range("K115").Value = "Hass"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 450
$ws.Range("N115").Value = 29000
$ws.Range("O115").Value = 30000
$ws.Range("P115").Value = 29444
$ws.Range("Q115").Value = "`$/bandeja 10 kilos"
$ws.Range("R115").Value = "Perú"
$ws.Range("S115").Value = 2944
$ws.Range("T115").Value = 10
$ws.Range("A116").Value = 1
$ws.Range("B116").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C116").Value = "Arica y Parinacota"
$ws.Range("D116").Value = 45016
$ws.Range("E116").Value = 15
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100106
$ws.Range("H116").Value = "Oleaginosos"
$ws.Range("I116").Value = 100106002
$ws.Range("J116").Value = "Palta"
$ws.Range("K116").Value = "Hass"
$ws.Range("L116").Value = "Segunda"
$ws.Range("M116").Value = 550
$ws.Range("N116").Value = 25000
$ws.Range("O116").Value = 26000
$ws.Range("P116").Value = 25364
$ws.Range("Q116").Value = "`$/bandeja 10 kilos"
$ws.Range("R116").Value = "Perú"
$ws.Range("S116").Value = 2536
$ws.Range("T116").Value = 10

Write-Output "Inserted 3 rows and populated new data."
